# Update vaccine calendar weekly delivery figures (Jan-Jul prognosis sheet).
# Each "master" row below drives a 7-day block via shared formulas in
# columns I:M (=D/7, =E/7, =F/7, =G/7, =H/7), which Excel recalculates
# automatically once the raw D/E/G/H inputs are edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week of 2021-04-05 (row 94) - Pfizer/Moderna split revised down
$ws.Range("D94").Value = 201240
$ws.Range("E94").Value = 34800

# Week of 2021-04-19 (row 108) - revised figures + number format applied
$ws.Range("D108").Value = 200070
$ws.Range("D108").NumberFormat = "#,##0"
$ws.Range("E108").Value = 22800

# Week of 2021-04-26 (row 115) - now also carries an E (Moderna) column
$ws.Range("D115").Value = 272610
$ws.Range("D115").NumberFormat = "#,##0"
$ws.Range("E115").Value = 22600

# Week of 2021-05-03 (row 122) - G/H (AZ columns) cleared, E added
$ws.Range("D122").Value = 267930
$ws.Range("E122").Value = 32400
$ws.Range("G122").ClearContents()
$ws.Range("H122").ClearContents()

# Week of 2021-05-10 (row 129) - H cleared, E added
$ws.Range("D129").Value = 273780
$ws.Range("E129").Value = 38400
$ws.Range("H129").ClearContents()

# Week of 2021-05-17 (row 136) - H cleared, E added
$ws.Range("D136").Value = 267930
$ws.Range("E136").Value = 32400
$ws.Range("H136").ClearContents()

# Week of 2021-05-24 (row 143) - H cleared, E added
$ws.Range("D143").Value = 274950
$ws.Range("E143").Value = 38400
$ws.Range("H143").ClearContents()

# Week of 2021-05-31 (row 150) - H cleared, E added
$ws.Range("D150").Value = 382590
$ws.Range("E150").Value = 31400
$ws.Range("H150").ClearContents()

# Week of 2021-06-07 (row 157) - new D/E inputs, H revised down
$ws.Range("D157").Value = 388440
$ws.Range("E157").Value = 38400
$ws.Range("H157").Value = 138867

# Week of 2021-06-14 (row 164) - new D/E inputs, H revised down
$ws.Range("D164").Value = 382590
$ws.Range("E164").Value = 32300
$ws.Range("H164").Value = 138867

# Week of 2021-06-21 (row 171) - new D/E inputs, H revised down
$ws.Range("D171").Value = 389610
$ws.Range("E171").Value = 44000
$ws.Range("H171").Value = 138867

# Week of 2021-06-28 (row 178) - new D/E inputs, H revised down
$ws.Range("D178").Value = 383760
$ws.Range("E178").Value = 36900
$ws.Range("H178").Value = 138867

# Weeks of 2021-07-05 / 07-12 / 07-19 / 07-26 (rows 185, 192, 199, 206) -
# H (AZ) revised down, D/E unchanged (blank)
$ws.Range("H185").Value = 683383
$ws.Range("H192").Value = 683383
$ws.Range("H199").Value = 683383
$ws.Range("H206").Value = 683383

# Restore the author's final cursor position / selection on the sheet
$ws.Range("F185").Select()
